# Auto-generated edit script applying scheduled price-update diff
# to Pandaemonium_Profits workbook (8 job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1092.2
$ws.Range("I19").Value = 832.1667
$ws.Range("K19").Value = 832.1667
$ws.Range("M19").Value = -657.1667
# Row 132
$ws.Range("H132").Value = 1162.4026
$ws.Range("I132").Value = 1196.4926
$ws.Range("J132").Value = 934
$ws.Range("K132").Value = 3589.4778
$ws.Range("L132").Value = 2802
$ws.Range("M132").Value = -1059.4778
$ws.Range("N132").Value = -7862
# Row 135
$ws.Range("H135").Value = 223.66667
$ws.Range("I135").Value = 201.48572
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 1813.37148
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = 721.6285200000002
$ws.Range("N135").Value = -14070
# Row 137
$ws.Range("H137").Value = 1525.9
$ws.Range("I137").Value = 1114.4822
$ws.Range("J137").Value = 3171.5715
$ws.Range("K137").Value = 3343.4466
$ws.Range("L137").Value = 9514.7145
$ws.Range("M137").Value = -793.4465999999998
$ws.Range("N137").Value = -14614.7145
# Row 138
$ws.Range("H138").Value = 3841.1526
$ws.Range("I138").Value = 2225.5715
$ws.Range("J138").Value = 4733.9736
$ws.Range("K138").Value = 6676.7145
$ws.Range("L138").Value = 14201.9208
$ws.Range("M138").Value = -1536.7145
$ws.Range("N138").Value = -24481.9208

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1846.909
$ws.Range("I45").Value = 1826.1333
$ws.Range("J45").Value = 1891.4286
$ws.Range("K45").Value = 1826.1333
$ws.Range("L45").Value = 1891.4286
$ws.Range("M45").Value = -1449.1333
$ws.Range("N45").Value = -2645.4286
# Row 74
$ws.Range("H74").Value = 4715.016
$ws.Range("I74").Value = 3295.3208
$ws.Range("K74").Value = 3295.3208
$ws.Range("M74").Value = -2421.3208
# Row 77
$ws.Range("H77").Value = 4715.016
$ws.Range("I77").Value = 3295.3208
$ws.Range("K77").Value = 16476.604
$ws.Range("M77").Value = -12108.604
# Row 97
$ws.Range("H97").Value = 1351.5264
$ws.Range("I97").Value = 1192.6364
$ws.Range("J97").Value = 1570
$ws.Range("K97").Value = 1192.6364
$ws.Range("L97").Value = 1570
$ws.Range("M97").Value = -696.6364000000001
$ws.Range("N97").Value = -2562
# Row 122
$ws.Range("H122").Value = 2841558.8
$ws.Range("I122").Value = 599.5854
$ws.Range("K122").Value = 1798.7562
$ws.Range("M122").Value = 651.2437999999997
# Row 132
$ws.Range("H132").Value = 3755.0566
$ws.Range("I132").Value = 1250.85
$ws.Range("J132").Value = 11460.308
$ws.Range("K132").Value = 3752.55
$ws.Range("L132").Value = 34380.924
$ws.Range("M132").Value = -1222.55
$ws.Range("N132").Value = -39440.924

$ws = $wb.Worksheets.Item("BSM")
# Row 40
$ws.Range("H40").Value = 40000
$ws.Range("J40").Value = 40000
$ws.Range("L40").Value = 40000
$ws.Range("N40").Value = -40530
# Row 94
$ws.Range("H94").Value = 1585.2
$ws.Range("I94").Value = 1601.125
$ws.Range("J94").Value = 1567
$ws.Range("K94").Value = 1601.125
$ws.Range("L94").Value = 1567
$ws.Range("M94").Value = -1150.125
$ws.Range("N94").Value = -2469
# Row 96
$ws.Range("H96").Value = 25950
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 25950
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 25950
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -31442
# Row 99
$ws.Range("H99").Value = 1291.7931
$ws.Range("I99").Value = 1092.6
$ws.Range("J99").Value = 1734.4445
$ws.Range("K99").Value = 1092.6
$ws.Range("L99").Value = 1734.4445
$ws.Range("M99").Value = 405.4000000000001
$ws.Range("N99").Value = -4730.4445

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3470.9355
$ws.Range("I31").Value = 3351.8542
$ws.Range("J31").Value = 3879.2144
$ws.Range("K31").Value = 3351.8542
$ws.Range("L31").Value = 3879.2144
$ws.Range("M31").Value = -3056.8542
$ws.Range("N31").Value = -4469.2144
# Row 34
$ws.Range("H34").Value = 3470.9355
$ws.Range("I34").Value = 3351.8542
$ws.Range("J34").Value = 3879.2144
$ws.Range("K34").Value = 3351.8542
$ws.Range("L34").Value = 3879.2144
$ws.Range("M34").Value = -3149.8542
$ws.Range("N34").Value = -4283.2144
# Row 58
$ws.Range("H58").Value = 1740.0597
$ws.Range("I58").Value = 1033.9387
$ws.Range("J58").Value = 3662.2778
$ws.Range("K58").Value = 1033.9387
$ws.Range("L58").Value = 3662.2778
$ws.Range("M58").Value = -830.9386999999999
$ws.Range("N58").Value = -4068.2778
# Row 86
$ws.Range("H86").Value = 2696.2964
$ws.Range("I86").Value = 2888.8823
$ws.Range("J86").Value = 2368.9
$ws.Range("K86").Value = 2888.8823
$ws.Range("L86").Value = 2368.9
$ws.Range("M86").Value = -1765.8823
$ws.Range("N86").Value = -4614.9
# Row 89
$ws.Range("H89").Value = 2696.2964
$ws.Range("I89").Value = 2888.8823
$ws.Range("J89").Value = 2368.9
$ws.Range("K89").Value = 14444.4115
$ws.Range("L89").Value = 11844.5
$ws.Range("M89").Value = -8828.411500000002
$ws.Range("N89").Value = -23076.5
# Row 132
$ws.Range("H132").Value = 6346.5
$ws.Range("I132").Value = 10398.385
$ws.Range("J132").Value = 3248
$ws.Range("K132").Value = 31195.155
$ws.Range("L132").Value = 9744
$ws.Range("M132").Value = -28665.155
$ws.Range("N132").Value = -14804
# Row 134
$ws.Range("H134").Value = 3008.0195
$ws.Range("I134").Value = 2208.6191
$ws.Range("J134").Value = 3567.6
$ws.Range("K134").Value = 6625.8573
$ws.Range("L134").Value = 10702.8
$ws.Range("M134").Value = -4090.8573
$ws.Range("N134").Value = -15772.8
# Row 136
$ws.Range("H136").Value = 1740.0597
$ws.Range("I136").Value = 1033.9387
$ws.Range("J136").Value = 3662.2778
$ws.Range("K136").Value = 3101.8161
$ws.Range("L136").Value = 10986.8334
$ws.Range("M136").Value = -551.8161
$ws.Range("N136").Value = -16086.8334

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 35714508
$ws.Range("I12").Value = 83333510
$ws.Range("J12").Value = 254.75
$ws.Range("K12").Value = 250000530
$ws.Range("L12").Value = 764.25
$ws.Range("M12").Value = -250000357
$ws.Range("N12").Value = -1110.25
# Row 122
$ws.Range("H122").Value = 834.5
$ws.Range("I122").Value = 623
$ws.Range("J122").Value = 985.5714
$ws.Range("K122").Value = 5607
$ws.Range("L122").Value = 8870.142600000001
$ws.Range("M122").Value = -3157
$ws.Range("N122").Value = -13770.1426

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 5011.7646
$ws.Range("I132").Value = 2982.3157
$ws.Range("J132").Value = 7582.4
$ws.Range("K132").Value = 8946.947100000001
$ws.Range("L132").Value = 22747.2
$ws.Range("M132").Value = -6416.947100000001
$ws.Range("N132").Value = -27807.2

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 828.5714
$ws.Range("I46").Value = 425
$ws.Range("J46").Value = 1366.6666
$ws.Range("K46").Value = 425
$ws.Range("L46").Value = 1366.6666
$ws.Range("M46").Value = -237
$ws.Range("N46").Value = -1742.6666
# Row 76
$ws.Range("H76").Value = 25747.445
$ws.Range("J76").Value = 25747.445
$ws.Range("L76").Value = 25747.445
$ws.Range("N76").Value = -26423.445
# Row 79
$ws.Range("H79").Value = 25747.445
$ws.Range("J79").Value = 25747.445
$ws.Range("L79").Value = 25747.445
$ws.Range("N79").Value = -28087.445
# Row 100
$ws.Range("H100").Value = 4058.3157
$ws.Range("I100").Value = 2282.182
$ws.Range("J100").Value = 6500.5
$ws.Range("K100").Value = 2282.182
$ws.Range("L100").Value = 6500.5
$ws.Range("M100").Value = -1741.182
$ws.Range("N100").Value = -7582.5
# Row 132
$ws.Range("H132").Value = 3289.403
$ws.Range("I132").Value = 3235.4814
$ws.Range("J132").Value = 3513.3845
$ws.Range("K132").Value = 9706.4442
$ws.Range("L132").Value = 10540.1535
$ws.Range("M132").Value = -7176.4442
$ws.Range("N132").Value = -15600.1535
# Row 136
$ws.Range("H136").Value = 4196.5454
$ws.Range("I136").Value = 2411.4102
$ws.Range("K136").Value = 7234.230599999999
$ws.Range("M136").Value = -4684.230599999999

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 34673.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 34673.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 34673.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -35303.5
# Row 73
$ws.Range("H73").Value = 34673.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 34673.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 34673.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -36857.5
# Row 132
$ws.Range("H132").Value = 1720.5278
$ws.Range("I132").Value = 855.5833
$ws.Range("J132").Value = 3450.4167
$ws.Range("K132").Value = 2566.7499
$ws.Range("L132").Value = 10351.2501
$ws.Range("M132").Value = -36.7498999999998
$ws.Range("N132").Value = -15411.2501

